# Refresh the cryptocurrency price/volume figures in the "Price" (D) and
# "Volume(1h)" (E) columns to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preserving its original
# number format/style (prevents Excel auto-detecting numeric-looking
# strings like "313.57" or "1.001" as actual numbers).
function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '27.284.55'
Set-TextValue 'E2' '  +0.55%  '
Set-TextValue 'D3' '1.774.83'
Set-TextValue 'E3' '  +3.78%  '
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '313.57'
Set-TextValue 'E6' '  +0.14%  '
Set-TextValue 'D7' '0.5245'
Set-TextValue 'E7' '  +9.62%  '
Set-TextValue 'D8' '0.3646'
Set-TextValue 'E8' '  +6.11%  '
Set-TextValue 'D9' '42.69'
Set-TextValue 'E9' '  +1.90%  '
Set-TextValue 'D10' '0.07365'
Set-TextValue 'E10' '  +0.98%  '
Set-TextValue 'E11' '  +4.13%  '
Set-TextValue 'E12' '  +0.09%  '
Set-TextValue 'D13' '20.48'
Set-TextValue 'E13' '  +3.13%  '
Set-TextValue 'E14' '  +3.66%  '
Set-TextValue 'D15' '1.766.66'
Set-TextValue 'E15' '  +3.58%  '
Set-TextValue 'D16' '6.953'
Set-TextValue 'E16' '  +1.82%  '
Set-TextValue 'D17' '88.70'
Set-TextValue 'E17' '  -0.54%  '
Set-TextValue 'E18' '  +0.44%  '
Set-TextValue 'D19' '0.06423'
Set-TextValue 'E19' '  +1.19%  '
Set-TextValue 'D20' '1.000'
Set-TextValue 'E20' '  +0.11%  '
Set-TextValue 'E21' '  +1.82%  '
Set-TextValue 'D22' '5.825'
Set-TextValue 'E22' '  +4.13%  '
Set-TextValue 'D23' '27.364.13'
Set-TextValue 'E23' '  +0.70%  '
Set-TextValue 'D24' '11.26'
Set-TextValue 'E24' '  +3.81%  '
Set-TextValue 'D25' '2.097'
Set-TextValue 'E25' '  -0.10%  '
Set-TextValue 'D26' '154.35'
Set-TextValue 'E26' '  -0.27%  '
Set-TextValue 'D27' '20.10'
Set-TextValue 'E27' '  +2.38%  '
Set-TextValue 'D28' '2.334'
Set-TextValue 'E28' '  +12.26%  '
Set-TextValue 'D29' '1.970.67'
Set-TextValue 'E29' '  +3.22%  '
Set-TextValue 'D30' '120.88'
Set-TextValue 'E30' '  +1.30%  '
Set-TextValue 'D31' '1.061'
Set-TextValue 'E31' '  +5.14%  '
Set-TextValue 'D32' '0.09762'
Set-TextValue 'E32' '  +5.57%  '
Set-TextValue 'D33' '5.555'
Set-TextValue 'E33' '  +4.70%  '
Set-TextValue 'D34' '3.621'
Set-TextValue 'E34' '  +1.24%  '
Set-TextValue 'D35' '0.02228'
Set-TextValue 'E35' '  +1.58%  '
Set-TextValue 'D36' '0.05960'
Set-TextValue 'E36' '  +1.83%  '
Set-TextValue 'D37' '11.22'
Set-TextValue 'E37' '  +1.50%  '
Set-TextValue 'D38' '4.836'
Set-TextValue 'E38' '  +2.01%  '
Set-TextValue 'E39' '  +1.48%  '
Set-TextValue 'E40' '  +4.33%  '
Set-TextValue 'D41' '1.433'
Set-TextValue 'E41' '  +2.19%  '
Set-TextValue 'D42' '8.061'
Set-TextValue 'E42' '  +8.24%  '
Set-TextValue 'D43' '1.140'
Set-TextValue 'E43' '  +2.88%  '
Set-TextValue 'D44' '13.13'
Set-TextValue 'E44' '  +4.78%  '
Set-TextValue 'D45' '3.626'
Set-TextValue 'E45' '  +1.96%  '
Set-TextValue 'D46' '0.5745'
Set-TextValue 'E46' '  +2.35%  '
Set-TextValue 'D47' '120.82'
Set-TextValue 'E47' '  +2.74%  '
Set-TextValue 'D48' '1.887'
Set-TextValue 'E48' '  +2.52%  '
Set-TextValue 'E49' '  +2.35%  '
Set-TextValue 'D50' '0.06713'
Set-TextValue 'E50' '  +1.32%  '
Set-TextValue 'D51' '70.47'
Set-TextValue 'E51' '  +1.07%  '
